# Applies the cryptos.xlsx price/volume update described in the commit diff
# (coinranking.com scrape refresh, GitHub Actions cron).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.422.21'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '2.248.03'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  +0.97%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("E9").Value = '  +1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.38'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0811'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").Value = '2.414.80'
$ws.Range("E14").Value = '  +6.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.840'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("D17").Value = '44.179.10'
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D18").Value = '0.0₃0967'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '66.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("E24").Value = '  +5.31%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  +6.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0800'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("E34").Value = '  -2.49%  '
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("E37").Value = '  +4.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.46'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.13%  '
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.23%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.85'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0303'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").Value = '1.753.75'
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.194'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '80.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '71.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.06%  '
$ws.Range("E50").Value = '  +7.20%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.00%  '

Write-Host "Applied cryptos list update"
